# Update the data table in Sheet1 with the new rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Index, B=Path, C=Symbol, D=Name, E=Email, F=Capital,
#          G=Optimization, H=Status, I=Change, J=TimeStamp
$rows = @(
    @{ Symbol="CRYPTO"; Name="Lorenzo Reyes"; Email="lreyes@udesa.edu.ar"; Capital="3200";    Optimization="MonteVaR";    TimeStamp="2022-10-18" },
    @{ Symbol="CRYPTO"; Name="Lorenzo Reyes"; Email="lreyes@udesa.edu.ar"; Capital="200000";  Optimization="MonteVaR";    TimeStamp="2022-10-08" },
    @{ Symbol="CRYPTO"; Name="Lorenzo Reyes"; Email="lreyes@udesa.edu.ar"; Capital="2985";    Optimization="MonteVaR";    TimeStamp="2022-10-08" },
    @{ Symbol="CRYPTO"; Name="Lorenzo Reyes"; Email="lreyes@udesa.edu.ar"; Capital="3200";    Optimization="MonteSharpe"; TimeStamp="2022-11-05" },
    @{ Symbol="CRYPTO"; Name="Peko Muni";     Email="pedrolandriel@live.com.ar"; Capital="1000000"; Optimization="SharpeRatio"; TimeStamp="2022-10-14" }
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $path = "./DATABASE/" + $row.Symbol + " " + $row.Name + " " + $row.Email + " " + $row.Capital + " " + $row.Optimization + " " + $row.TimeStamp + ".xlsx"

    $ws.Cells.Item($r, 2).Value = $path
    $ws.Cells.Item($r, 3).Value = $row.Symbol
    $ws.Cells.Item($r, 4).Value = $row.Name
    $ws.Cells.Item($r, 5).Value = $row.Email

    # Columns F (Capital) and J (TimeStamp) hold numeric-/date-looking text
    # (e.g. "3200", "2022-10-18") that must stay plain text rather than be
    # auto-converted to a number or date, so enter them with a leading
    # apostrophe (Excel's normal "force text" entry method).
    $ws.Cells.Item($r, 6).Value = "'" + $row.Capital
    $ws.Cells.Item($r, 7).Value = $row.Optimization
    $ws.Cells.Item($r, 10).Value = "'" + $row.TimeStamp
}
